$d = $word.ActiveDocument

# The static "www.drpaulduenas.com" text in the footer needs to become a
# MERGEFIELD-driven "website" field (begin/instrText/separate/result/end),
# mirroring the other MERGEFIELD-based runs already used in this footer
# (e.g. =emergency_number, =consultation.branch_office.address, etc.)
# All of the new runs keep the exact same run formatting the literal text
# run used to have.

$oldText = "www.drpaulduenas.com"

$targetFooter = $null
foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists -and ($ftr.Range.Text -like "*$oldText*")) {
            $targetFooter = $ftr
            break
        }
    }
    if ($targetFooter -ne $null) { break }
}

if ($targetFooter -eq $null) {
    throw "Could not locate a footer containing '$oldText'"
}

$find = $targetFooter.Range.Find
# Replace:=0 (wdReplaceNone) -- locate only, do not let Execute perform its
# own replacement; we need the matched Range intact so InsertXML can target it.
$found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '$oldText' to replace with the website field"
}

$matched = $find.Parent

$rPr = '<w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'

$runBegin     = '<w:r>' + $rPr + '<w:fldChar w:fldCharType="begin"/></w:r>'
$runInstr     = '<w:r>' + $rPr + '<w:instrText xml:space="preserve"> MERGEFIELD =website \* MERGEFORMAT </w:instrText></w:r>'
$runSeparate  = '<w:r>' + $rPr + '<w:fldChar w:fldCharType="separate"/></w:r>'
$runResult    = '<w:r>' + $rPr + '<w:t>' + [char]0x00AB + '=website' + [char]0x00BB + '</w:t></w:r>'
$runEnd       = '<w:r>' + $rPr + '<w:fldChar w:fldCharType="end"/></w:r>'

$bodyRuns = $runBegin + $runInstr + $runSeparate + $runResult + $runEnd

$xmlFragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p>' + $bodyRuns + '</w:p></w:body>' + `
    '</w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

# Insert the new field runs right next to the matched ("www.drpaulduenas.com")
# range, preserving the paragraph they live in, then remove the old literal
# text run that triggered the match.
$matched.InsertXML($xmlFragment)

$cleanup = $targetFooter.Range.Find
$cleaned = $cleanup.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
if (-not $cleaned) {
    throw "Failed to remove the original '$oldText' run after inserting the website field"
}

Write-Output "website field inserted: $($targetFooter.Range.Text)"
